$d = $word.ActiveDocument

# The signed date on the Independent Completion Form changes from
# "08/29/2021" to "09/04/2021". The day/month portion ("08/29") is
# replaced with "09/04"; the "/2021" that follows is left untouched,
# matching the target edit exactly.
$d.Content.Find.Execute("08/29", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "09/04", 2)
